$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Week 1 lecture row: swap the old "Lecture 0" overview content
# for the new "Lecture 1" introduction-to-linear-models content.
$ws.Range("C2").Value = "Lecture 1:  Introduction to Linear Models"
$ws.Range("D2").Value = "01-reading.html"
$ws.Range("E2").Value = "01-introduction"

# The reading link that used to live on row 3 (01-reading.html) has moved
# up to row 2, so clear it out of its old spot.
$ws.Range("D3").Value = ""

# Update the lab link for the first lab session.
$ws.Range("F4").Value = "lab-01.html"

# Restore the selection to the last-edited cell.
$ws.Range("F4").Select() | Out-Null
